# 190331_Research Ideas for Topics of Social Science.docx
# 1) Title paragraph: append " Term Paper " (bold, Raleway, 12pt/sz24) after
#    "...Social Science".
# 2) Second paragraph: replace the "Deadline for project proposal..." text
#    with "Group members: ..." (italic, Raleway, 12pt/sz24), and move the
#    hidden "_GoBack" bookmark here (it previously sat further down, around
#    "the spread of ").

$d = $word.ActiveDocument

# --- 1) Title paragraph -----------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.Find.Execute(
    "Research Ideas for Topics of Social Science", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Research Ideas for Topics of Social Science Term Paper ", 2) | Out-Null

# --- 2) Second paragraph -------------------------------------------------
$deadlinePara = $d.Paragraphs.Item(2)
$deadlinePara.Range.Find.Execute(
    "Deadline for project proposal: 26th of April 2019", $true, $false,
    $false, $false, $false, $true, 1, $false,
    "Group members: Sina Smid, Edith Zink, Zeyu Zhao, Helge Zille", 2) | Out-Null

# Re-fetch the paragraph (its Range grew) and stamp sz/szCs=24 (12pt) onto
# the paragraph mark and every run inside it, matching the rest of the
# formatting already present (Raleway italic, en-US).
$deadlinePara = $d.Paragraphs.Item(2)
$deadlineRange = $deadlinePara.Range
$deadlineRange.Font.Size = 12
$deadlineRange.Font.SizeBi = 12

# Move the "_GoBack" bookmark: delete it from its old location further down
# the document and re-create it spanning the new paragraph's text.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$deadlinePara = $d.Paragraphs.Item(2)
$deadlineRange = $deadlinePara.Range
$bookmarkRange = $d.Range($deadlineRange.Start, $deadlineRange.End - 1)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null
